$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A51").Style = "source"
Write-Host "Font.Italic: " $ws.Range("A51").Font.Italic
Write-Host "Font.Underline: " $ws.Range("A51").Font.Underline
